# Add the missing city ("ErnieTown") for Ernie's driver row, completing the
# table of drivers (row 3, column G = "city").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "ErnieTown"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("G4").Select()
